$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# Adds a new handed-back file, f1fd48bb-ecd5-45a8-aa2b-2572d2da1df2.md, as a
# new row (row 4) on each of the three worksheets:
#   - "Overview" (7-column summary table)
#   - "zh-cn"     (16-column detail table)
#   - "de-de"     (16-column detail table)
# ---------------------------------------------------------------------------

$baseId      = "f1fd48bb-ecd5-45a8-aa2b-2572d2da1df2"
$mdName      = "$baseId.md"
$mdPath      = "e2e\$baseId.md"
$statusSync  = "Handed back: in sync with en-US"
$xlfHash     = "33abd1804a59ff02c3173677d87a96fe90d3179f"
$zhcnXlf     = "$baseId.$xlfHash.zh-cn.xlf"
$dedeXlf     = "$baseId.$xlfHash.de-de.xlf"

$zhcnHandoffDt  = "2016-08-29 20:58:17"
$zhcnHandbackDt = "2016-08-29 20:58:34"
$dedeHandoffDt  = "2016-08-29 20:58:22"
$dedeHandbackDt = "2016-08-29 20:58:42"

$overviewCommit = "97abdbe760bfb783d5d673011df14962b4404870"
$zhcnCommit     = "ca201ec60e5cd5c1be6d6e6a70c8065f90b80245"
$dedeCommit     = "df99ef42c7917f0d51cd99df5fa8abb7d4961f49"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# 1) Overview sheet
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$rowOv = $loOv.ListRows.Add()
$rngOv = $rowOv.Range

$rngOv.Cells.Item(1,1).Value = $mdName
$bCell = $rngOv.Cells.Item(1,2)
$bCell.Value = $mdPath
$hOv = $wsOv.Hyperlinks.Add($bCell, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$overviewCommit/e2e/$mdName", "", "", $mdPath)
$rngOv.Cells.Item(1,3).Value = ".md"
$rngOv.Cells.Item(1,5).Value = $statusSync
$rngOv.Cells.Item(1,6).Value = $statusSync
$gCell = $rngOv.Cells.Item(1,7)
$gCell.NumberFormat = $dateFmt
$gCell.Value = $dedeHandoffDt

# ---------------------------------------------------------------------------
# 2) zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()
$rngZh = $rowZh.Range

$aCellZh = $rngZh.Cells.Item(1,1)
$aCellZh.Value = $mdName
$hZhA = $wsZh.Hyperlinks.Add($aCellZh, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$overviewCommit/$mdPath", "", "", $mdName)

$rngZh.Cells.Item(1,2).Value = ".md"
$rngZh.Cells.Item(1,3).Value = $statusSync
$rngZh.Cells.Item(1,4).Value = "e2e"
$rngZh.Cells.Item(1,5).Value = "ht"
$rngZh.Cells.Item(1,6).Value = "'True"
$rngZh.Cells.Item(1,7).Value = $zhcnXlf
$hCellZh = $rngZh.Cells.Item(1,8)
$hCellZh.NumberFormat = $dateFmt
$hCellZh.Value = $zhcnHandoffDt

$iCellZh = $rngZh.Cells.Item(1,9)
$iCellZh.Value = $mdName
$hZhI = $wsZh.Hyperlinks.Add($iCellZh, "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$zhcnCommit/e2e/$mdName", "", "", $mdName)

$rngZh.Cells.Item(1,10).Value = $zhcnXlf
$kCellZh = $rngZh.Cells.Item(1,11)
$kCellZh.NumberFormat = $dateFmt
$kCellZh.Value = $zhcnHandbackDt

$rngZh.Cells.Item(1,12).Value = "'"
$rngZh.Cells.Item(1,13).Value = "'True"
$rngZh.Cells.Item(1,14).Value = "'"
$rngZh.Cells.Item(1,15).Value = "'False"
$rngZh.Cells.Item(1,16).Value = "'"

# ---------------------------------------------------------------------------
# 3) de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()
$rngDe = $rowDe.Range

$aCellDe = $rngDe.Cells.Item(1,1)
$aCellDe.Value = $mdName
$hDeA = $wsDe.Hyperlinks.Add($aCellDe, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$overviewCommit/$mdPath", "", "", $mdName)

$rngDe.Cells.Item(1,2).Value = ".md"
$rngDe.Cells.Item(1,3).Value = $statusSync
$rngDe.Cells.Item(1,4).Value = "e2e"
$rngDe.Cells.Item(1,5).Value = "ht"
$rngDe.Cells.Item(1,6).Value = "'True"
$rngDe.Cells.Item(1,7).Value = $dedeXlf
$hCellDe = $rngDe.Cells.Item(1,8)
$hCellDe.NumberFormat = $dateFmt
$hCellDe.Value = $dedeHandoffDt

$iCellDe = $rngDe.Cells.Item(1,9)
$iCellDe.Value = $mdName
$hDeI = $wsDe.Hyperlinks.Add($iCellDe, "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$dedeCommit/e2e/$mdName", "", "", $mdName)

$rngDe.Cells.Item(1,10).Value = $dedeXlf
$kCellDe = $rngDe.Cells.Item(1,11)
$kCellDe.NumberFormat = $dateFmt
$kCellDe.Value = $dedeHandbackDt

$rngDe.Cells.Item(1,12).Value = "'"
$rngDe.Cells.Item(1,13).Value = "'True"
$rngDe.Cells.Item(1,14).Value = "'"
$rngDe.Cells.Item(1,15).Value = "'False"
$rngDe.Cells.Item(1,16).Value = "'"

Write-Output "Handback row added for $mdName"
